# Apply crypto price/volume update (GitHub Actions automated refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need an explicit Text
# number format first, otherwise Excel will silently convert the string
# into a numeric value (and drop formatting such as trailing zeros).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Assign the updated cell text
$ws.Range("D2").Value = "59.179.76"
$ws.Range("E2").Value = "  -2.87%  "
$ws.Range("D3").Value = "2.650.55"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "523.03"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").Value = "144.85"
$ws.Range("E6").Value = "  -2.56%  "
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("E8").Value = "  -1.31%  "
$ws.Range("D9").Value = "7.02"
$ws.Range("E9").Value = "  +9.23%  "
$ws.Range("E10").Value = "  -3.78%  "
$ws.Range("E11").Value = "  -2.27%  "
$ws.Range("E12").Value = "  +1.83%  "
$ws.Range("D13").Value = "3.119.97"
$ws.Range("E13").Value = "  -1.38%  "
$ws.Range("D14").Value = "59.217.03"
$ws.Range("D15").Value = "21.17"
$ws.Range("E15").Value = "  -1.41%  "
$ws.Range("E16").Value = "  -2.39%  "
$ws.Range("D17").Value = "2.657.40"
$ws.Range("E17").Value = "  -7.63%  "
$ws.Range("D18").Value = "339.98"
$ws.Range("E18").Value = "  -4.31%  "
$ws.Range("E19").Value = "  -4.29%  "
$ws.Range("D20").Value = "10.37"
$ws.Range("E20").Value = "  -1.95%  "
$ws.Range("D21").Value = "6.36"
$ws.Range("E21").Value = "  +0.68%  "
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").Value = "63.78"
$ws.Range("E23").Value = "  +1.76%  "
$ws.Range("E24").Value = "  -2.78%  "
$ws.Range("D25").Value = "0.167"
$ws.Range("E25").Value = "  -1.35%  "
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  +0.65%  "
$ws.Range("D27").Value = "0.0₃0802"
$ws.Range("E27").Value = "  -3.12%  "
$ws.Range("D28").Value = "7.10"
$ws.Range("E28").Value = "  -2.58%  "
$ws.Range("E29").Value = "  -2.37%  "
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("D32").Value = "18.82"
$ws.Range("E32").Value = "  -1.73%  "
$ws.Range("D33").Value = "149.11"
$ws.Range("E33").Value = "  -0.78%  "
$ws.Range("D34").Value = "4.16"
$ws.Range("E34").Value = "  -1.24%  "
$ws.Range("E35").Value = "  -2.07%  "
$ws.Range("D36").Value = "0.897"
$ws.Range("E36").Value = "  -6.08%  "
$ws.Range("D37").Value = "0.880"
$ws.Range("E37").Value = "  -0.50%  "
$ws.Range("D38").Value = "36.73"
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("D39").Value = "1.48"
$ws.Range("E39").Value = "  -6.05%  "
$ws.Range("E40").Value = "  -3.98%  "
$ws.Range("D41").Value = "0.617"
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.51%  "
$ws.Range("D43").Value = "275.34"
$ws.Range("E43").Value = "  -3.97%  "
$ws.Range("D44").Value = "19.84"
$ws.Range("E44").Value = "  -0.81%  "
$ws.Range("D45").Value = "0.0972"
$ws.Range("E45").Value = "  -2.17%  "
$ws.Range("E47").Value = "  +2.15%  "
$ws.Range("D48").Value = "4.79"
$ws.Range("E48").Value = "  -2.99%  "
$ws.Range("D49").Value = "2.028.62"
$ws.Range("E49").Value = "  -5.64%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "0.0228"
$ws.Range("E50").Value = "  -2.81%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "18.90"
$ws.Range("E51").Value = "  -2.29%  "
